$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("Z1").Value = "t1"
$f1 = $ws.Range("Z1").Font
$f1.Name = "Arial"
$f1.Size = 14
$f1.ThemeColor = 1
$f1.TintAndShade = 0.14999847407452621

$ws.Range("Z2").Value = "t2"
$f2 = $ws.Range("Z2").Font
$f2.Name = "Arial"
$f2.Size = 10
$f2.ThemeColor = 9
$f2.TintAndShade = -0.499984740745262
Write-Host "done"
